$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'36.497.73"
$ws.Range("E2").Value = "  +1.26%  "
$ws.Range("D3").Value = "'1.943.51"
$ws.Range("E3").Value = "  -1.12%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'243.59"
$ws.Range("E5").Value = "  +0.31%  "
$ws.Range("D6").Value = "'0.614"
$ws.Range("E6").Value = "  -1.10%  "
$ws.Range("D7").Value = "'57.89"
$ws.Range("E7").Value = "  -5.31%  "
$ws.Range("E8").Value = "  -0.14%  "
$ws.Range("D9").Value = "'0.368"
$ws.Range("E9").Value = "  -0.53%  "
$ws.Range("D10").Value = "'55.71"
$ws.Range("E10").Value = "  -0.85%  "
$ws.Range("D11").Value = "'0.0840"
$ws.Range("E11").Value = "  +6.21%  "
$ws.Range("E12").Value = "  +1.05%  "
$ws.Range("D13").Value = "'21.68"
$ws.Range("E13").Value = "  -0.55%  "
$ws.Range("D14").Value = "'0.825"
$ws.Range("E14").Value = "  -3.81%  "
$ws.Range("D15").Value = "'2.230.63"
$ws.Range("E15").Value = "  -1.41%  "
$ws.Range("D16").Value = "'13.60"
$ws.Range("E16").Value = "  -1.97%  "
$ws.Range("D17").Value = "'5.26"
$ws.Range("E17").Value = "  -2.29%  "
$ws.Range("D18").Value = "'1.942.20"
$ws.Range("E18").Value = "  -1.95%  "
$ws.Range("D19").Value = "'36.430.16"
$ws.Range("E19").Value = "  +1.28%  "
$ws.Range("B20").Value = "Litecoin"
$ws.Range("C20").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D20").Value = "'69.69"
$ws.Range("E20").Value = "  -1.38%  "
$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").Value = "'0.0₃0870"
$ws.Range("E21").Value = "  +2.54%  "
$ws.Range("D22").Value = "'229.88"
$ws.Range("E22").Value = "  -3.05%  "
$ws.Range("D23").Value = "'5.02"
$ws.Range("E23").Value = "  -3.12%  "
$ws.Range("E24").Value = "  +0.12%  "
$ws.Range("D25").Value = "'2.45"
$ws.Range("E25").Value = "  -1.34%  "
$ws.Range("E26").Value = "  +1.32%  "
$ws.Range("D27").Value = "'9.27"
$ws.Range("E27").Value = "  -3.97%  "
$ws.Range("D28").Value = "'162.27"
$ws.Range("E28").Value = "  +2.40%  "
$ws.Range("D29").Value = "'19.42"
$ws.Range("E29").Value = "  -0.96%  "
$ws.Range("E30").Value = "  -7.59%  "
$ws.Range("E31").Value = "  -0.63%  "
$ws.Range("D32").Value = "'1.15"
$ws.Range("E32").Value = "  +1.92%  "
$ws.Range("D33").Value = "'4.68"
$ws.Range("E33").Value = "  -3.40%  "
$ws.Range("D34").Value = "'0.0630"
$ws.Range("E34").Value = "  +2.86%  "
$ws.Range("D35").Value = "'4.28"
$ws.Range("E35").Value = "  -1.44%  "
$ws.Range("D36").Value = "'6.26"
$ws.Range("E36").Value = "  +0.77%  "
$ws.Range("E37").Value = "  +0.08%  "
$ws.Range("E38").Value = "  -3.31%  "
$ws.Range("D39").Value = "'2.14"
$ws.Range("E39").Value = "  -5.68%  "
$ws.Range("D40").Value = "'3.05"
$ws.Range("E40").Value = "  -1.25%  "
$ws.Range("D41").Value = "'0.0972"
$ws.Range("E41").Value = "  -0.65%  "
$ws.Range("E42").Value = "  +4.81%  "
$ws.Range("E43").Value = "  -3.89%  "
$ws.Range("E44").Value = "  -1.37%  "
$ws.Range("D45").Value = "'16.07"
$ws.Range("E45").Value = "  +0.06%  "
$ws.Range("D46").Value = "'1.354.59"
$ws.Range("E46").Value = "  +1.47%  "
$ws.Range("D47").Value = "'1.03"
$ws.Range("E47").Value = "  -4.57%  "
$ws.Range("D48").Value = "'87.70"
$ws.Range("E48").Value = "  -4.60%  "
$ws.Range("D49").Value = "'7.16"
$ws.Range("E49").Value = "  -4.25%  "
$ws.Range("E50").Value = "  +2.48%  "
$ws.Range("D51").Value = "'45.16"
$ws.Range("E51").Value = "  +3.66%  "
